$wb = $excel.ActiveWorkbook

# --- Sheet "main": remove item names from B1:H1, keep A1 only ---
$wsMain = $wb.Worksheets.Item("main")
$wsMain.Range("B1:H1").ClearContents()
$wsMain.Range("B1").Select() | Out-Null

# --- Sheet "items": add a new column header + two new rows ---
$wsItems = $wb.Worksheets.Item("items")
$wsItems.Range("D1").Value = "배고픔 회복량"
$wsItems.Range("A14").Value = "카드키"
$wsItems.Range("B14").Value = 10
$wsItems.Range("C14").Value = "보안실, 창고, 관리실 등에서 사용 가능"
$wsItems.Range("A15").Value = "빵"
$wsItems.Range("B15").Value = 2
$wsItems.Range("C15").Value = "배고픔 20 회복"
$wsItems.Range("D15").Value = 20
$wsItems.Range("E16").Select() | Out-Null

# --- Sheet "b_hall": add an item to box 2 ---
$wsBHall = $wb.Worksheets.Item("b_hall")
$wsBHall.Range("B4").Value = "빵"
$wsBHall.Range("E5").Select() | Out-Null

# --- Sheet "sp3": no content change, selection stays B2 ---
$wsSp3 = $wb.Worksheets.Item("sp3")
$wsSp3.Range("B2").Select() | Out-Null

# b_hall becomes the active tab (activeTab index 2 / 3rd sheet)
$wsBHall.Activate() | Out-Null
